$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data to match the latest GitHub Actions scrape.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '23.569.40'
$ws.Range('E2').Value = '  +1.66%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.657.74'
$ws.Range('E3').Value = '  +2.88%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.34%  '
$ws.Range('E5').Value = '  -0.15%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '302.27'
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3838'
$ws.Range('E7').Value = '  +1.52%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '51.07'
$ws.Range('E8').Value = '  -1.33%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3594'
$ws.Range('E9').Value = '  +1.82%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.242'
$ws.Range('E10').Value = '  +3.41%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08191'
$ws.Range('E11').Value = '  +1.28%  '
$ws.Range('E12').Value = '  -0.21%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.42'
$ws.Range('E13').Value = '  +1.83%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.503'
$ws.Range('E14').Value = '  +2.22%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.503'
$ws.Range('E15').Value = '  +3.80%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001225'
$ws.Range('E16').Value = '  +1.53%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.654.27'
$ws.Range('E17').Value = '  +2.16%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '97.46'
$ws.Range('E18').Value = '  +3.43%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06973'
$ws.Range('E19').Value = '  +0.89%  '
$ws.Range('E20').Value = '  +5.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.71'
$ws.Range('E21').Value = '  +3.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').Value = '  -0.09%  '
$ws.Range('E23').Value = '  +3.33%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '23.566.89'
$ws.Range('E24').Value = '  +1.66%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.500'
$ws.Range('E25').Value = '  -0.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.996'
$ws.Range('E26').Value = '  -0.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.21'
$ws.Range('E27').Value = '  +1.80%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '151.83'
$ws.Range('E28').Value = '  +0.67%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.233'
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.95'
$ws.Range('E30').Value = '  +1.29%  '
$ws.Range('B31').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C31').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.835.66'
$ws.Range('E31').Value = '  +2.01%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.197'
$ws.Range('E32').Value = '  +11.40%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.244'
$ws.Range('E33').Value = '  +6.90%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.98'
$ws.Range('E34').Value = '  +4.80%  '
$ws.Range('E35').Value = '  -1.68%  '
$ws.Range('E36').Value = '  +3.71%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.135'
$ws.Range('E37').Value = '  +5.36%  '
$ws.Range('B38').Value = 'Algorand'
$ws.Range('C38').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2498'
$ws.Range('E38').Value = '  +2.13%  '
$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.08793'
$ws.Range('E39').Value = '  +0.69%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.07032'
$ws.Range('E40').Value = '  +1.65%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '13.17'
$ws.Range('E41').Value = '  +10.36%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.7032'
$ws.Range('E42').Value = '  +2.64%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.338'
$ws.Range('E43').Value = '  +1.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.11'
$ws.Range('E44').Value = '  +5.18%  '
$ws.Range('E45').Value = '  +4.26%  '
$ws.Range('E46').Value = '  -0.10%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.309'
$ws.Range('E47').Value = '  +3.20%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.958'
$ws.Range('E48').Value = '  +0.45%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.07937'
$ws.Range('E49').Value = '  +1.15%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '127.01'
$ws.Range('E50').Value = '  +0.16%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.194'
$ws.Range('E51').Value = '  +2.49%  '
